$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-10-08 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-09 Monday", 2)

$d.Content.Find.Execute("80÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "65÷4=", 2)
$d.Content.Find.Execute("14÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=", 2)
$d.Content.Find.Execute("61÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷6=", 2)
$d.Content.Find.Execute("34÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷4=", 2)
$d.Content.Find.Execute("28÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷4=", 2)

$d.Content.Find.Execute("59÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷6=", 2)
$d.Content.Find.Execute("89÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷6=", 2)
$d.Content.Find.Execute("40÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=", 2)
$d.Content.Find.Execute("92÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷7=", 2)
$d.Content.Find.Execute("81÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷4=", 2)

$d.Content.Find.Execute("70÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷9=", 2)
$d.Content.Find.Execute("82÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷2=", 2)
$d.Content.Find.Execute("95÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷4=", 2)
$d.Content.Find.Execute("80÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷4=", 2)
$d.Content.Find.Execute("17÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷9=", 2)

$d.Content.Find.Execute("33÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=", 2)
$d.Content.Find.Execute("59÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷8=", 2)
$d.Content.Find.Execute("76÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷2=", 2)
$d.Content.Find.Execute("60÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷8=", 2)
$d.Content.Find.Execute("24÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "43÷6=", 2)

$d.Content.Find.Execute("13÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷7=", 2)
$d.Content.Find.Execute("26÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷8=", 2)
$d.Content.Find.Execute("42÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷7=", 2)
$d.Content.Find.Execute("87÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷6=", 2)
$d.Content.Find.Execute("84÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=", 2)
